# The "BLS Data Series" sheet holds the monthly unemployment-rate matrix,
# one row per year (2000-2020). The source data series does not actually
# start until 2003, so the rows for 2000, 2001 and 2002 are removed,
# shifting everything else up by three rows.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BLS Data Series")
$ws2 = $wb.Worksheets.Item("Sheet1")

# Rows 2-4 hold 2000, 2001, 2002 (row 1 is the Jan..Dec header).
$ws1.Range("A2:M4").EntireRow.Delete()

# Restore Sheet1's own selection before switching away from it, so it
# doesn't keep a dangling "active" selection once BLS Data Series takes
# over as the active/selected tab.
$ws2.Range("A1:XFD12").Select()

# Make the data sheet the active tab, with the (now shifted) first data
# rows selected, matching the post-edit author view.
$ws1.Activate()
$ws1.Range("A2:XFD4").Select()
